$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('CAN_Canada', 90),
    @('KEN_Kenya', 90),
    @('FRA_France', 90),
    @('CHE_Switzerland', 90),
    @('PHL_Philippines', 90),
    @('JPN_Japan', 90),
    @('COL_Colombia', 90),
    @('ITA_Italy', 90),
    @('ISR_Israel', 90),
    @('CRI_Costa Rica', 90),
    @('SWE_Sweden', 90),
    @('PER_Peru', 90),
    @('ISL_Iceland', 90),
    @('DEU_Germany', 90),
    @('IRL_Ireland', 90),
    @('IND_India', 90),
    @('DNK_Denmark', 90),
    @('PRT_Portugal', 90),
    @('ECU_Ecuador', 90),
    @('EGY_Egypt', 90),
    @('ESP_Spain', 90),
    @('GTM_Guatemala', 90),
    @('FIN_Finland', 90),
    @('CYP_Cyprus', 90),
    @('THA_Thailand', 90),
    @('GBR_United Kingdom', 90),
    @('BRA_Brazil', 90),
    @('MUS_Mauritius', 90),
    @('ZAF_South Africa', 90),
    @('VEN_Venezuela (Bolivarian Republic of)', 90),
    @('ARG_Argentina', 90),
    @('USA_United States', 90),
    @('URY_Uruguay', 90),
    @('AUS_Australia', 90),
    @('AUT_Austria', 90),
    @('MEX_Mexico', 90),
    @('NLD_Netherlands', 90),
    @('MAR_Morocco', 90),
    @('NOR_Norway', 90),
    @('BEL_Belgium', 90),
    @('LKA_Sri Lanka', 90),
    @('BOL_Bolivia (Plurinational State of)', 90),
    @('NZL_New Zealand', 90),
    @('TTO_Trinidad and Tobago', 90),
    @('TUR_Turkey', 90),
    @('LUX_Luxembourg', 90),
    @('CHL_Chile', 102),
    @('GRC_Greece', 108),
    @('PRY_Paraguay', 108),
    @('DOM_Dominican Republic', 108),
    @('TWN_Taiwan', 108),
    @('CHN_China', 126),
    @('KOR_Republic of Korea', 144),
    @('JAM_Jamaica', 144),
    @('NGA_Nigeria', 156),
    @('MLT_Malta', 162),
    @('JOR_Jordan', 162),
    @('ZWE_Zimbabwe', 162),
    @('IRN_Iran (Islamic Republic of)', 180),
    @('PAN_Panama', 180),
    @('ZMB_Zambia', 180),
    @('MYS_Malaysia', 180),
    @('HND_Honduras', 186),
    @('IDN_Indonesia', 246),
    @('NIC_Nicaragua', 246),
    @('BFA_Burkina Faso', 255),
    @('BWA_Botswana', 270),
    @('BRB_Barbados', 270),
    @('TUN_Tunisia', 270),
    @('SGP_Singapore', 270),
    @('SEN_Senegal', 270),
    @('TZA_U.R. of Tanzania: Mainland', 270),
    @('CIV_Côte d''Ivoire', 270),
    @('GAB_Gabon', 270),
    @('NER_Niger', 270),
    @('MOZ_Mozambique', 270),
    @('ROU_Romania', 270),
    @('RWA_Rwanda', 270),
    @('NAM_Namibia', 270),
    @('HKG_China, Hong Kong SAR', 270),
    @('CMR_Cameroon', 270),
    @('MRT_Mauritania', 348),
    @('BEN_Benin', 354),
    @('LSO_Lesotho', 366),
    @('BDI_Burundi', 366),
    @('CAF_Central African Republic', 366),
    @('TGO_Togo', 366),
    @('FJI_Fiji', 366),
    @('SLE_Sierra Leone', 378),
    @('POL_Poland', 450),
    @('SDN_Sudan', 450),
    @('BHR_Bahrain', 450),
    @('HUN_Hungary', 450),
    @('BGR_Bulgaria', 450),
    @('SAU_Saudi Arabia', 450),
    @('IRQ_Iraq', 450),
    @('QAT_Qatar', 450),
    @('AGO_Angola', 450),
    @('KWT_Kuwait', 450),
    @('PAK_Pakistan', 464),
    @('UGA_Uganda', 464),
    @('SLV_El Salvador', 464),
    @('COD_D.R. of the Congo', 464),
    @('MNG_Mongolia', 486),
    @('MAC_China, Macao SAR', 486),
    @('SWZ_Eswatini', 486),
    @('LAO_Lao People''s DR', 486),
    @('ETH_Ethiopia', 494),
    @('MWI_Malawi', 512),
    @('GHA_Ghana', 524),
    @('BGD_Bangladesh', 572),
    @('NPL_Nepal', 584),
    @('SYR_Syrian Arab Republic', 584),
    @('GMB_Gambia', 584),
    @('COG_Congo', 584),
    @('MLI_Mali', 584),
    @('HTI_Haiti', 584),
    @('MDG_Madagascar', 584),
    @('DZA_Algeria', 584),
    @('MMR_Myanmar', 608),
    @('LBR_Liberia', 632),
    @('VNM_Viet Nam', 704),
    @('KHM_Cambodia', 704),
    @('MDV_Maldives', 704),
    @('BRN_Brunei Darussalam', 704),
    @('BLZ_Belize', 704),
    @('ALB_Albania', 704),
    @('ARE_United Arab Emirates', 704),
    @('GUY_Guyana', 704),
    @('GIN_Guinea', 744),
    @('GNB_Guinea-Bissau', 753),
    @('CPV_Cabo Verde', 753),
    @('TCD_Chad', 753),
    @('SYC_Seychelles', 753),
    @('GNQ_Equatorial Guinea', 753),
    @('COM_Comoros', 753),
    @('KAZ_Kazakhstan', 810),
    @('MDA_Republic of Moldova', 810),
    @('SVN_Slovenia', 810),
    @('LVA_Latvia', 810),
    @('TJK_Tajikistan', 810),
    @('SRB_Serbia', 810),
    @('HRV_Croatia', 810),
    @('LTU_Lithuania', 810),
    @('CZE_Czech Republic', 810),
    @('UKR_Ukraine', 810),
    @('RUS_Russian Federation', 810),
    @('KGZ_Kyrgyzstan', 810),
    @('ARM_Armenia', 810),
    @('EST_Estonia', 810),
    @('SVK_Slovakia', 810),
    @('BHS_Bahamas', 843),
    @('VGB_British Virgin Islands', 843),
    @('BMU_Bermuda', 843),
    @('VCT_St. Vincent and the Grenadines', 843),
    @('BTN_Bhutan', 843),
    @('ATG_Antigua and Barbuda', 843),
    @('AIA_Anguilla', 843),
    @('TCA_Turks and Caicos Islands', 843),
    @('ABW_Aruba', 843),
    @('MSR_Montserrat', 843),
    @('SUR_Suriname', 843),
    @('STP_Sao Tome and Principe', 843),
    @('CYM_Cayman Islands', 843),
    @('DJI_Djibouti', 843),
    @('LCA_Saint Lucia', 843),
    @('DMA_Dominica', 843),
    @('OMN_Oman', 843),
    @('LBN_Lebanon', 843),
    @('PSE_State of Palestine', 843),
    @('GRD_Grenada', 843),
    @('KNA_Saint Kitts and Nevis', 843),
    @('YEM_Yemen', 932),
    @('GEO_Georgia', 1023),
    @('MKD_North Macedonia', 1023),
    @('MNE_Montenegro', 1023),
    @('AZE_Azerbaijan', 1023),
    @('BIH_Bosnia and Herzegovina', 1023),
    @('BLR_Belarus', 1023),
    @('TKM_Turkmenistan', 1023),
    @('UZB_Uzbekistan', 1023),
    @('SXM_Sint Maarten (Dutch part)', 1158),
    @('CUW_Curaçao', 1158),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}